# Export Map Overlay MAX/MIN values into the MaxMin worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# INNER MAX / MIN (feet) values
$ws.Range("C5").Value = 710.12
$ws.Range("C6").Value = 603.058

# OUTER MAX / MIN (feet) values
$ws.Range("C11").Value = 726.69
$ws.Range("C12").Value = 599.527

$wb.Save()
